$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3045
$ws.Range("I64").Value = 2813.8572
$ws.Range("J64").Value = 3449.5
$ws.Range("K64").Value = 2813.8572
$ws.Range("L64").Value = 3449.5
$ws.Range("M64").Value = -2565.8572
$ws.Range("N64").Value = -3945.5
$ws.Range("H67").Value = 3045
$ws.Range("I67").Value = 2813.8572
$ws.Range("J67").Value = 3449.5
$ws.Range("K67").Value = 2813.8572
$ws.Range("L67").Value = 3449.5
$ws.Range("M67").Value = -1955.8572
$ws.Range("N67").Value = -5165.5
$ws.Range("H76").Value = 10422654
$ws.Range("I76").Value = 6637.7407
$ws.Range("K76").Value = 6637.7407
$ws.Range("M76").Value = -6322.7407
$ws.Range("H79").Value = 10422654
$ws.Range("I79").Value = 6637.7407
$ws.Range("K79").Value = 6637.7407
$ws.Range("M79").Value = -5545.7407

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 450
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 450
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -682
$ws.Range("H61").Value = 1906.2858
$ws.Range("I61").Value = 2084.4443
$ws.Range("J61").Value = 1802.8387
$ws.Range("K61").Value = 2084.4443
$ws.Range("L61").Value = 1802.8387
$ws.Range("M61").Value = -1872.4443
$ws.Range("N61").Value = -2226.8387
$ws.Range("H63").Value = 2110.25
$ws.Range("I63").Value = 2110.25
$ws.Range("K63").Value = 2110.25
$ws.Range("M63").Value = -1424.25
$ws.Range("H66").Value = 2110.25
$ws.Range("I66").Value = 2110.25
$ws.Range("K66").Value = 10551.25
$ws.Range("M66").Value = -7119.25
$ws.Range("H74").Value = 2095.639
$ws.Range("I74").Value = 2518.375
$ws.Range("J74").Value = 1250.1666
$ws.Range("K74").Value = 2518.375
$ws.Range("L74").Value = 1250.1666
$ws.Range("M74").Value = -1644.375
$ws.Range("N74").Value = -2998.1666
$ws.Range("H77").Value = 2095.639
$ws.Range("I77").Value = 2518.375
$ws.Range("J77").Value = 1250.1666
$ws.Range("K77").Value = 12591.875
$ws.Range("L77").Value = 6250.833000000001
$ws.Range("M77").Value = -8223.875
$ws.Range("N77").Value = -14986.833
$ws.Range("H97").Value = 1126.5
$ws.Range("I97").Value = 1023.5
$ws.Range("J97").Value = 1216.625
$ws.Range("K97").Value = 1023.5
$ws.Range("L97").Value = 1216.625
$ws.Range("M97").Value = -527.5
$ws.Range("N97").Value = -2208.625
$ws.Range("H132").Value = 6214.9565
$ws.Range("I132").Value = 1528.8572
$ws.Range("J132").Value = 13504.444
$ws.Range("K132").Value = 4586.571599999999
$ws.Range("L132").Value = 40513.33199999999
$ws.Range("M132").Value = -2056.571599999999
$ws.Range("N132").Value = -45573.33199999999
$ws.Range("H136").Value = 1906.2858
$ws.Range("I136").Value = 2084.4443
$ws.Range("J136").Value = 1802.8387
$ws.Range("K136").Value = 6253.3329
$ws.Range("L136").Value = 5408.5161
$ws.Range("M136").Value = -3703.3329
$ws.Range("N136").Value = -10508.5161

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4577.154
$ws.Range("I20").Value = 4466.222
$ws.Range("K20").Value = 4466.222
$ws.Range("M20").Value = -4219.222
$ws.Range("H105").Value = 31252378
$ws.Range("I105").Value = 2427.12
$ws.Range("J105").Value = 142859340
$ws.Range("K105").Value = 2427.12
$ws.Range("L105").Value = 142859340
$ws.Range("M105").Value = -680.1199999999999
$ws.Range("N105").Value = -142862834

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2561.6086
$ws.Range("I31").Value = 1693.325
$ws.Range("J31").Value = 3759.2415
$ws.Range("K31").Value = 1693.325
$ws.Range("L31").Value = 3759.2415
$ws.Range("M31").Value = -1398.325
$ws.Range("N31").Value = -4349.2415
$ws.Range("H34").Value = 2561.6086
$ws.Range("I34").Value = 1693.325
$ws.Range("J34").Value = 3759.2415
$ws.Range("K34").Value = 1693.325
$ws.Range("L34").Value = 3759.2415
$ws.Range("M34").Value = -1491.325
$ws.Range("N34").Value = -4163.2415
$ws.Range("H62").Value = 3579.4075
$ws.Range("I62").Value = 2832.7144
$ws.Range("J62").Value = 4383.5386
$ws.Range("K62").Value = 2832.7144
$ws.Range("L62").Value = 4383.5386
$ws.Range("M62").Value = -2208.7144
$ws.Range("N62").Value = -5631.5386
$ws.Range("H65").Value = 3579.4075
$ws.Range("I65").Value = 2832.7144
$ws.Range("J65").Value = 4383.5386
$ws.Range("K65").Value = 14163.572
$ws.Range("L65").Value = 21917.693
$ws.Range("M65").Value = -11043.572
$ws.Range("N65").Value = -28157.693

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 723.2
$ws.Range("I131").Value = 417.27274
$ws.Range("J131").Value = 761.0112
$ws.Range("K131").Value = 1251.81822
$ws.Range("L131").Value = 2283.0336
$ws.Range("M131").Value = 3788.18178
$ws.Range("N131").Value = -12363.0336

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4044
$ws.Range("I70").Value = 4125.579
$ws.Range("J70").Value = 3914.8333
$ws.Range("K70").Value = 4125.579
$ws.Range("L70").Value = 3914.8333
$ws.Range("M70").Value = -3855.579
$ws.Range("N70").Value = -4454.8333
$ws.Range("H73").Value = 4044
$ws.Range("I73").Value = 4125.579
$ws.Range("J73").Value = 3914.8333
$ws.Range("K73").Value = 4125.579
$ws.Range("L73").Value = 3914.8333
$ws.Range("M73").Value = -3189.579
$ws.Range("N73").Value = -5786.8333
$ws.Range("H80").Value = 7695972.5
$ws.Range("I80").Value = 3970.3333
$ws.Range("K80").Value = 3970.3333
$ws.Range("M80").Value = -2972.3333
$ws.Range("H83").Value = 7695972.5
$ws.Range("I83").Value = 3970.3333
$ws.Range("K83").Value = 19851.6665
$ws.Range("M83").Value = -14859.6665

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1337.1724
$ws.Range("I61").Value = 1243.3914
$ws.Range("J61").Value = 1696.6666
$ws.Range("K61").Value = 1243.3914
$ws.Range("L61").Value = 1696.6666
$ws.Range("M61").Value = -1041.3914
$ws.Range("N61").Value = -2100.6666
$ws.Range("H113").Value = 1337.1724
$ws.Range("I113").Value = 1243.3914
$ws.Range("J113").Value = 1696.6666
$ws.Range("K113").Value = 1243.3914
$ws.Range("L113").Value = 1696.6666
$ws.Range("M113").Value = 926.6086
$ws.Range("N113").Value = -6036.6666
$ws.Range("H136").Value = 3949.077
$ws.Range("I136").Value = 2864.6072
$ws.Range("J136").Value = 6709.5454
$ws.Range("K136").Value = 8593.821599999999
$ws.Range("L136").Value = 20128.6362
$ws.Range("M136").Value = -6043.821599999999
$ws.Range("N136").Value = -25228.6362

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14724118
$ws.Range("I132").Value = 18538936
$ws.Range("J132").Value = 9814.429
$ws.Range("K132").Value = 55616808
$ws.Range("L132").Value = 29443.287
$ws.Range("M132").Value = -55614278
$ws.Range("N132").Value = -34503.287
$ws.Range("H136").Value = 3368.2354
$ws.Range("I136").Value = 4448.067
$ws.Range("J136").Value = 1825.619
$ws.Range("K136").Value = 13344.201
$ws.Range("L136").Value = 5476.857
$ws.Range("M136").Value = -10794.201
$ws.Range("N136").Value = -10576.857
